$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: new period labels Q1:S1 ---
$ws.Range("Q1").Value = "31/12/2023"
$ws.Range("R1").Value = "31/03/2024"
$ws.Range("S1").Value = "30/06/2024"
$ws.Range("P1").Copy()
$ws.Range("Q1:S1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Numeric data rows: Q/R/S values ---
$ws.Range("Q2").Value = 351209.984
$ws.Range("R2").Value = 333081.984
$ws.Range("S2").Value = 319304
$ws.Range("Q3").Value = 211022
$ws.Range("R3").Value = 193726
$ws.Range("S3").Value = 185216
$ws.Range("Q4").Value = 139348
$ws.Range("R4").Value = 131916
$ws.Range("S4").Value = 129261
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("Q6").Value = 44328
$ws.Range("R6").Value = 39494
$ws.Range("S6").Value = 37294
$ws.Range("Q7").Value = 22879
$ws.Range("R7").Value = 20320
$ws.Range("S7").Value = 17425
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("Q9").Value = 3096
$ws.Range("R9").Value = 468
$ws.Range("S9").Value = 310
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("Q11").Value = 1371
$ws.Range("R11").Value = 1528
$ws.Range("S11").Value = 926
$ws.Range("Q12").Value = 91702
$ws.Range("R12").Value = 94459
$ws.Range("S12").Value = 93473
$ws.Range("Q13").Value = 2649
$ws.Range("R13").Value = 2578
$ws.Range("S13").Value = 2503
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("S16").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 0
$ws.Range("S18").Value = 0
$ws.Range("Q19").Value = 54889
$ws.Range("R19").Value = 55680
$ws.Range("S19").Value = 56386
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("S20").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = 0
$ws.Range("S21").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("R22").Value = 0
$ws.Range("S22").Value = 0
$ws.Range("Q23").Value = 32919
$ws.Range("R23").Value = 29768
$ws.Range("S23").Value = 25938
$ws.Range("Q24").Value = 15567
$ws.Range("R24").Value = 15129
$ws.Range("S24").Value = 14677
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = 0
$ws.Range("S25").Value = 0
$ws.Range("Q26").Value = 351209.984
$ws.Range("R26").Value = 333081.984
$ws.Range("S26").Value = 319304
$ws.Range("Q27").Value = 73863
$ws.Range("R27").Value = 67324
$ws.Range("S27").Value = 63141
$ws.Range("Q28").Value = 9400
$ws.Range("R28").Value = 6736
$ws.Range("S28").Value = 6898
$ws.Range("Q29").Value = 39349
$ws.Range("R29").Value = 32309
$ws.Range("S29").Value = 29918
$ws.Range("Q30").Value = 3662
$ws.Range("R30").Value = 3357
$ws.Range("S30").Value = 2720
$ws.Range("Q31").Value = 9202
$ws.Range("R31").Value = 9904
$ws.Range("S31").Value = 10122
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = 0
$ws.Range("S32").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = 0
$ws.Range("S33").Value = 0
$ws.Range("Q34").Value = 12250
$ws.Range("R34").Value = 15018
$ws.Range("S34").Value = 13483
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = 0
$ws.Range("S35").Value = 0
$ws.Range("Q36").Value = 0
$ws.Range("R36").Value = 0
$ws.Range("S36").Value = 0
$ws.Range("Q37").Value = 21943
$ws.Range("R37").Value = 24858
$ws.Range("S37").Value = 23525
$ws.Range("Q38").Value = 6565
$ws.Range("R38").Value = 5655
$ws.Range("S38").Value = 3506
$ws.Range("Q39").Value = 0
$ws.Range("R39").Value = 0
$ws.Range("S39").Value = 0
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = 4940
$ws.Range("S40").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = 0
$ws.Range("S41").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = 0
$ws.Range("S42").Value = 0
$ws.Range("Q43").Value = 15378
$ws.Range("R43").Value = 14263
$ws.Range("S43").Value = 20019
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = 0
$ws.Range("S44").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("R45").Value = 0
$ws.Range("S45").Value = 0
$ws.Range("Q46").Value = 0
$ws.Range("R46").Value = 0
$ws.Range("S46").Value = 0
$ws.Range("Q47").Value = 255404
$ws.Range("R47").Value = 240900
$ws.Range("S47").Value = 232638
$ws.Range("Q48").Value = 471166.016
$ws.Range("R48").Value = 471230.016
$ws.Range("S48").Value = 471284
$ws.Range("Q49").Value = -10431
$ws.Range("R49").Value = -10404
$ws.Range("S49").Value = -10380
$ws.Range("Q50").Value = 0
$ws.Range("R50").Value = 0
$ws.Range("S50").Value = 0
$ws.Range("Q51").Value = 0
$ws.Range("R51").Value = 0
$ws.Range("S51").Value = 0
$ws.Range("Q52").Value = -205331.008
$ws.Range("R52").Value = -219926
$ws.Range("S52").Value = -228266
$ws.Range("Q53").Value = 0
$ws.Range("R53").Value = 0
$ws.Range("S53").Value = 0
$ws.Range("Q54").Value = 0
$ws.Range("R54").Value = 0
$ws.Range("S54").Value = 0
$ws.Range("Q55").Value = 0
$ws.Range("R55").Value = 0
$ws.Range("S55").Value = 0
$ws.Range("Q56").Value = 0
$ws.Range("R56").Value = 0
$ws.Range("S56").Value = 0
$ws.Range("Q59").Value = 56948
$ws.Range("R59").Value = 40712
$ws.Range("S59").Value = 41504
$ws.Range("Q60").Value = -38040
$ws.Range("R60").Value = -24348
$ws.Range("S60").Value = -23667
$ws.Range("Q61").Value = 18908
$ws.Range("R61").Value = 16364
$ws.Range("S61").Value = 17837
$ws.Range("Q62").Value = -30014
$ws.Range("R62").Value = -14867
$ws.Range("S62").Value = -10091
$ws.Range("Q63").Value = -13000
$ws.Range("R63").Value = -22153
$ws.Range("S63").Value = -20648
$ws.Range("Q64").Value = 0
$ws.Range("R64").Value = 0
$ws.Range("S64").Value = 0
$ws.Range("Q65").Value = 3747
$ws.Range("R65").Value = 457
$ws.Range("S65").Value = 1491
$ws.Range("Q66").Value = 0
$ws.Range("R66").Value = 0
$ws.Range("S66").Value = 0
$ws.Range("Q67").Value = 0
$ws.Range("R67").Value = 0
$ws.Range("S67").Value = 0
$ws.Range("Q68").Value = 2809
$ws.Range("R68").Value = 4813
$ws.Range("S68").Value = 2365
$ws.Range("Q69").Value = 4379
$ws.Range("R69").Value = 5384
$ws.Range("S69").Value = 4045
$ws.Range("Q70").Value = -1570
$ws.Range("R70").Value = -571
$ws.Range("S70").Value = -1680
$ws.Range("Q74").Value = -17550
$ws.Range("R74").Value = -15386
$ws.Range("S74").Value = -9046
$ws.Range("Q75").Value = 0
$ws.Range("R75").Value = 0
$ws.Range("S75").Value = 0
$ws.Range("Q76").Value = 1015
$ws.Range("R76").Value = 0
$ws.Range("S76").Value = 706
$ws.Range("Q80").Value = -16535
$ws.Range("R80").Value = -14595
$ws.Range("S80").Value = -8340

# --- Blank separator rows: create formatted blank cells Q:S matching existing blank pattern ---
$ws.Range("B57:D57").Copy()
$ws.Range("Q57").PasteSpecial(-4122)
$ws.Range("B58:D58").Copy()
$ws.Range("Q58").PasteSpecial(-4122)
$ws.Range("B71:D71").Copy()
$ws.Range("Q71").PasteSpecial(-4122)
$ws.Range("B72:D72").Copy()
$ws.Range("Q72").PasteSpecial(-4122)
$ws.Range("B73:D73").Copy()
$ws.Range("Q73").PasteSpecial(-4122)
$ws.Range("B77:D77").Copy()
$ws.Range("Q77").PasteSpecial(-4122)
$ws.Range("B78:D78").Copy()
$ws.Range("Q78").PasteSpecial(-4122)
$ws.Range("B79:D79").Copy()
$ws.Range("Q79").PasteSpecial(-4122)
$excel.CutCopyMode = $false
